# "Add Leave Card 10/32023 3:18 PM"
# Adds newly-credited leave entries (Aug-Sep 2023) to the
# "2017 LEAVE BALANCE" sheet (rows 68-74) and the corresponding
# EARNED value on the "2018 LEAVE CREDITS" sheet (rows 83-85).
# Dependent formulas (BALANCE columns, CONVERTION!A7) recalc automatically.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) "2018 LEAVE CREDITS" sheet - EARNED column entries for the
#    three bi-monthly periods ending 8/1/2023, 9/1/2023, 10/2/2023
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("2018 LEAVE CREDITS")
$ws1.Range("C83").Value = 1.25
$ws1.Range("C84").Value = 1.25
$ws1.Range("C85").Value = 1.25
$ws1.Range("F90").Select()

# ---------------------------------------------------------------
# 2) "2017 LEAVE BALANCE" sheet - new SL/VL leave rows 68-74
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("2017 LEAVE BALANCE")

# Row 68 - SL(1-0-0) taken 8/15/2023, credited 8/5/2023
$ws2.Range("A68").Value = 45153
$ws2.Range("B68").Value = "SL(1-0-0)"
$ws2.Range("H68").Value = 1
$ws2.Range("K68").NumberFormat = "m/d/yy"
$ws2.Range("K68").Value = 45143

# Row 69 - VL(2-0-0) taken 8/12,15/2023
$ws2.Range("B69").Value = "VL(2-0-0)"
$ws2.Range("D69").Value = 2
$ws2.Range("K69").Value = "8/12,15/2023"

# Row 70 - SL(1-0-0) credited 8/13/2023
$ws2.Range("B70").Value = "SL(1-0-0)"
$ws2.Range("H70").Value = 1
$ws2.Range("K70").NumberFormat = "m/d/yy"
$ws2.Range("K70").Value = 45151

# Row 71 - VL(3-0-0) taken 8/22-24/2023
$ws2.Range("B71").Value = "VL(3-0-0)"
$ws2.Range("D71").Value = 3
$ws2.Range("K71").Value = "8/22-24/2023"

# Row 72 - SL(2-0-0) taken 8/31 , 9/2/2023
$ws2.Range("A72").Value = 45170
$ws2.Range("B72").Value = "SL(2-0-0)"
$ws2.Range("H72").Value = 2
$ws2.Range("K72").Value = "8/31 , 9/2/2023"

# Row 73 - VL(3-0-0) taken 9/27-29/2023
$ws2.Range("B73").Value = "VL(3-0-0)"
$ws2.Range("D73").Value = 3
$ws2.Range("K73").Value = "9/27-29/2023"

# Row 74 - VL(3-0-0) taken 9/19-21/2023
$ws2.Range("B74").Value = "VL(3-0-0)"
$ws2.Range("D74").Value = 3
$ws2.Range("K74").Value = "9/19-21/2023"

$ws2.Range("K74").Select()
